# Adds two new columns, I ("I0") and J ("IF"), with per-row numeric data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1/J1: copy H1 formatting (bold header w/ border) then set text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-68: column I = "I0", column J = "IF".
$ijData = @{
  2 = @(7, 8)
  3 = @(4, 5)
  4 = @(6, 6)
  5 = @(7, 7)
  6 = @(7, 7)
  7 = @(4, 5)
  8 = @(8, 8)
  9 = @(5, 6)
  10 = @(7, 7)
  11 = @(7, 7)
  12 = @(9, 9)
  13 = @(8, 8)
  14 = @(6, 6)
  15 = @(7, 7)
  16 = @(5, 6)
  17 = @(6, 7)
  18 = @(8, 8)
  19 = @(8, 8)
  20 = @(9, 9)
  21 = @(6, 7)
  22 = @(8, 8)
  23 = @(9, 9)
  24 = @(9, 9)
  25 = @(9, 9)
  26 = @(7, 8)
  27 = @(8, 8)
  28 = @(9, 9)
  29 = @(9, 9)
  30 = @(7, 7)
  31 = @(7, 8)
  32 = @(9, 9)
  33 = @(8, 8)
  34 = @(6, 7)
  35 = @(9, 10)
  36 = @(8, 8)
  37 = @(5, 6)
  38 = @(8, 8)
  39 = @(8, 8)
  40 = @(7, 7)
  41 = @(6, 7)
  42 = @(7, 8)
  43 = @(6, 7)
  44 = @(8, 8)
  45 = @(7, 8)
  46 = @(7, 7)
  47 = @(5, 6)
  48 = @(9, 9)
  49 = @(8, 8)
  50 = @(8, 8)
  51 = @(3, 5)
  52 = @(7, 7)
  53 = @(10, 11)
  54 = @(5, 5)
  55 = @(1, 1)
  56 = @(6, 8)
  57 = @(3, 3)
  58 = @(7, 7)
  59 = @(3, 4)
  60 = @(8, 9)
  61 = @(3, 5)
  62 = @(4, 6)
  63 = @(6, 6)
  64 = @(1, 3)
  65 = @(6, 7)
  66 = @(7, 7)
  67 = @(7, 8)
  68 = @(3, 4)
}

foreach ($row in $ijData.Keys) {
  $pair = $ijData[$row]
  $ws.Cells.Item($row, 9).Value = $pair[0]
  $ws.Cells.Item($row, 10).Value = $pair[1]
}
